# Apply the "Trade #182 closed" update to the live_trading_results workbook.
#
# Summary of changes:
#  - Summary sheet: Total Trades 181 -> 182, Win Rate % 41.99 -> 41.76
#  - Strategy Status sheet: volatility_scorer row -> Trades 8 -> 9, Win Rate % 50 -> 44.44
#  - All Trades sheet: append the two new trade rows (closed #182, opened #183)
#  - volatility_scorer sheet: append the newly closed trade row (#182)
#  - MarketMaking sheet: append the newly opened trade row (#183)

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($ws, $row, $col, $value)
    # Force a text number format first so date/time-looking strings
    # ("2026-02-17", "10:03:19") are not auto-converted into date/time serials.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 182
$wsSummary.Range("B9").Value = 41.76

# ---------------------------------------------------------------------------
# Strategy Status sheet (volatility_scorer row = row 12)
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D12").Value = 9
$wsStatus.Range("G12").Value = 44.44

# ---------------------------------------------------------------------------
# All Trades sheet - append row 183 (closed trade #182) and row 184 (opened trade #183)
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$r = 183
$wsAll.Cells.Item($r, 1).Value = 182
Set-TextCell $wsAll $r 2 "2026-02-17"
Set-TextCell $wsAll $r 3 "10:03:19"
$wsAll.Cells.Item($r, 4).Value = "volatility_scorer"
$wsAll.Cells.Item($r, 5).Value = "NEUTRAL"
$wsAll.Cells.Item($r, 6).Value = 0.97
$wsAll.Cells.Item($r, 7).Value = 0.97
$wsAll.Cells.Item($r, 8).Value = "CLOSED"
$wsAll.Cells.Item($r, 9).Value = 0
$wsAll.Cells.Item($r, 10).Value = 0
$wsAll.Cells.Item($r, 11).Value = 100.11
$wsAll.Cells.Item($r, 12).Value = 0
$wsAll.Cells.Item($r, 13).Value = 0
$wsAll.Cells.Item($r, 14).Value = 0.85
$wsAll.Cells.Item($r, 15).Value = "Low vol market (score: inf) - ideal for market making"
$wsAll.Cells.Item($r, 16).Value = "early_exit"
$wsAll.Cells.Item($r, 17).Value = 0.18

$r = 184
$wsAll.Cells.Item($r, 1).Value = 183
Set-TextCell $wsAll $r 2 "2026-02-17"
Set-TextCell $wsAll $r 3 "10:03:19"
$wsAll.Cells.Item($r, 4).Value = "MarketMaking"
$wsAll.Cells.Item($r, 5).Value = "DOWN"
$wsAll.Cells.Item($r, 6).Value = 0.97
$wsAll.Cells.Item($r, 8).Value = "OPEN"
$wsAll.Cells.Item($r, 9).Value = 0
$wsAll.Cells.Item($r, 10).Value = 0
$wsAll.Cells.Item($r, 11).Value = 100.6480687506789
$wsAll.Cells.Item($r, 12).Value = 0
$wsAll.Cells.Item($r, 13).Value = 0
$wsAll.Cells.Item($r, 14).Value = 0.6
$wsAll.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
$wsAll.Cells.Item($r, 17).Value = 0

# ---------------------------------------------------------------------------
# volatility_scorer sheet - append row 10 (closed trade #182, same data as above)
# ---------------------------------------------------------------------------
$wsVol = $wb.Worksheets.Item("volatility_scorer")

$r = 10
$wsVol.Cells.Item($r, 1).Value = 182
Set-TextCell $wsVol $r 2 "2026-02-17"
Set-TextCell $wsVol $r 3 "10:03:19"
$wsVol.Cells.Item($r, 4).Value = "volatility_scorer"
$wsVol.Cells.Item($r, 5).Value = "NEUTRAL"
$wsVol.Cells.Item($r, 6).Value = 0.97
$wsVol.Cells.Item($r, 7).Value = 0.97
$wsVol.Cells.Item($r, 8).Value = "CLOSED"
$wsVol.Cells.Item($r, 9).Value = 0
$wsVol.Cells.Item($r, 10).Value = 0
$wsVol.Cells.Item($r, 11).Value = 100.11
$wsVol.Cells.Item($r, 12).Value = 0
$wsVol.Cells.Item($r, 13).Value = 0
$wsVol.Cells.Item($r, 14).Value = 0.85
$wsVol.Cells.Item($r, 15).Value = "Low vol market (score: inf) - ideal for market making"
$wsVol.Cells.Item($r, 16).Value = "early_exit"
$wsVol.Cells.Item($r, 17).Value = 0.18

# ---------------------------------------------------------------------------
# MarketMaking sheet - append row 175 (opened trade #183, same data as above)
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

$r = 175
$wsMM.Cells.Item($r, 1).Value = 183
Set-TextCell $wsMM $r 2 "2026-02-17"
Set-TextCell $wsMM $r 3 "10:03:19"
$wsMM.Cells.Item($r, 4).Value = "MarketMaking"
$wsMM.Cells.Item($r, 5).Value = "DOWN"
$wsMM.Cells.Item($r, 6).Value = 0.97
$wsMM.Cells.Item($r, 8).Value = "OPEN"
$wsMM.Cells.Item($r, 9).Value = 0
$wsMM.Cells.Item($r, 10).Value = 0
$wsMM.Cells.Item($r, 11).Value = 100.6480687506789
$wsMM.Cells.Item($r, 12).Value = 0
$wsMM.Cells.Item($r, 13).Value = 0
$wsMM.Cells.Item($r, 14).Value = 0.6
$wsMM.Cells.Item($r, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item($r, 17).Value = 0

Write-Host "Applied trade #182/#183 update."
